# Royal Flush Description.docx — adds the "Pong Game" and "Hobbies"
# sections (with their body paragraphs) after the existing final
# paragraph, right before the section properties.

$d = $word.ActiveDocument

# Anchor on the very end of the document body (just before </w:body>,
# i.e. right before the sectPr). This is where the new paragraphs go.
$insertAt = $d.Range($d.Content.End, $d.Content.End)

# The new content is expressed as literal WordprocessingML (several
# paragraphs, including a couple of bold/size-28 headings, some blank
# paragraphs used as spacers, and a w:lastRenderedPageBreak marker in
# the final paragraph) and injected verbatim via Range.InsertXML so the
# resulting markup matches exactly rather than relying on the object
# model to synthesize equivalent-but-different XML.
$newParagraphsXml = @'
<w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Pong Game</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>After building this website, I realized that I had not gotten as much practice with JavaScript as I had hoped. I decided that I’d build a simple pong game as a means to get to know some of the functions, syntax, and techniques unique to JavaScript. The game I ended up building has progressive levels that get harder and harder along with a health system that keeps track of how many lives the player has left. I defiantly feel like this project did the job. I learned a lot and got a better understanding of what JavaScript can do!</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Hobbies</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>In the spring of 2018, I bought my first motorcycle. I fell in love with riding and I would ride it everywhere. When I rode it to school it would wake me up in the morning and get me ready to learn. When I rode my bike home all the stress of that day would melt away. After a few months, my natural inclination to fixing and improving things kicked in.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>My first project was to change the clutch of the bike. As the bike was already 20 years old the clutch had worn and had begun to slip. After watching instructional videos and reading the shop manual I began the work. This was the first major repair project I had ever done, and it took much longer than expected. Over the following months I would tinker with my bike on Sundays, culminating with me disassembling my friends’ bike that didn’t run and getting it running again. The image shows one of my projects where I changed out the cam chain tensioner, which sets the timing of the engine.</w:t></w:r></w:p><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>With my growing confidence, and love for working on physical projects, I would fix things around the house whenever needed. While studying abroad in Israel, my friends room flooded pretty bad and completely submerged my friends MacBook underwater.  After retrieving the laptop, he didn’t know what to do. Being in a foreign country he didn’t know who he could get to fix his water damaged laptop. Having nowhere else to turn he asked me to try and fix it. I told him I was almost positive I wouldn’t be able to, but I would give it a try. The process stared from the beginning, research on “how does water damage ruin electronics”, reading various articles and developing a game plan. I decided that the best I could do is to thoroughly clean off any corrosion on all of the connections with alcohol. That is what I did. After gingerly disassembling the laptop, cleaning is, and then putting it back together. I held my breath as I powered it up… IT WORKED!</w:t></w:r></w:p>
'@

$packageXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          "<w:body>$newParagraphsXml</w:body>" +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertAt.InsertXML($packageXml)

Write-Host "Paragraph count after insert:" $d.Paragraphs.Count
